$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 211.26471
$ws.Range("I33").Value = 168.21875
$ws.Range("J33").Value = 900
$ws.Range("K33").Value = 168.21875
$ws.Range("L33").Value = 900
$ws.Range("M33").Value = 60.78125
$ws.Range("N33").Value = -1358
$ws.Range("H57").Value = 28749.5
$ws.Range("J57").Value = 28749.5
$ws.Range("L57").Value = 86248.5
$ws.Range("N57").Value = -87246.5
$ws.Range("H98").Value = 2091.6
$ws.Range("I98").Value = 2232.8108
$ws.Range("J98").Value = 350
$ws.Range("K98").Value = 2232.8108
$ws.Range("L98").Value = 350
$ws.Range("M98").Value = -734.8108000000002
$ws.Range("N98").Value = -3346
$ws.Range("H112").Value = 30304420
$ws.Range("I112").Value = 606
$ws.Range("J112").Value = 35715816
$ws.Range("K112").Value = 1818
$ws.Range("L112").Value = 107147448
$ws.Range("M112").Value = -710
$ws.Range("N112").Value = -107149664
$ws.Range("H121").Value = 1000
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -6494
$ws.Range("H122").Value = 2091.6
$ws.Range("I122").Value = 2232.8108
$ws.Range("J122").Value = 350
$ws.Range("K122").Value = 6698.432400000001
$ws.Range("L122").Value = 1050
$ws.Range("M122").Value = -4248.432400000001
$ws.Range("N122").Value = -5950
$ws.Range("H131").Value = 1676.75
$ws.Range("I131").Value = 1133.125
$ws.Range("J131").Value = 3851.25
$ws.Range("K131").Value = 3399.375
$ws.Range("L131").Value = 11553.75
$ws.Range("M131").Value = 1640.625
$ws.Range("N131").Value = -21633.75
$ws.Range("H135").Value = 1195.2709
$ws.Range("I135").Value = 282.1389
$ws.Range("J135").Value = 3934.6667
$ws.Range("K135").Value = 2539.2501
$ws.Range("L135").Value = 35412.0003
$ws.Range("M135").Value = -4.250099999999748
$ws.Range("N135").Value = -40482.0003
$ws.Range("H137").Value = 33334780
$ws.Range("I137").Value = 1368.5385
$ws.Range("J137").Value = 250001950
$ws.Range("K137").Value = 4105.6155
$ws.Range("L137").Value = 750005850
$ws.Range("M137").Value = -1555.6155
$ws.Range("N137").Value = -750010950
$ws.Range("H141").Value = 627.37933
$ws.Range("I141").Value = 578.3571
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 1735.0713
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 3444.9287
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4687.169
$ws.Range("I32").Value = 4966.544
$ws.Range("K32").Value = 4966.544
$ws.Range("M32").Value = -4679.544
$ws.Range("H74").Value = 3819.8718
$ws.Range("I74").Value = 843.8276
$ws.Range("J74").Value = 12450.4
$ws.Range("K74").Value = 843.8276
$ws.Range("L74").Value = 12450.4
$ws.Range("M74").Value = 30.17240000000004
$ws.Range("N74").Value = -14198.4
$ws.Range("H77").Value = 3819.8718
$ws.Range("I77").Value = 843.8276
$ws.Range("J77").Value = 12450.4
$ws.Range("K77").Value = 4219.138
$ws.Range("L77").Value = 62252
$ws.Range("M77").Value = 148.8620000000001
$ws.Range("N77").Value = -70988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2043
$ws.Range("I86").Value = 1765.7142
$ws.Range("J86").Value = 2572.3635
$ws.Range("K86").Value = 1765.7142
$ws.Range("L86").Value = 2572.3635
$ws.Range("M86").Value = -642.7141999999999
$ws.Range("N86").Value = -4818.363499999999
$ws.Range("H89").Value = 2043
$ws.Range("I89").Value = 1765.7142
$ws.Range("J89").Value = 2572.3635
$ws.Range("K89").Value = 8828.571
$ws.Range("L89").Value = 12861.8175
$ws.Range("M89").Value = -3212.571
$ws.Range("N89").Value = -24093.8175

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2463.5
$ws.Range("I134").Value = 2624.9412
$ws.Range("J134").Value = 2071.4285
$ws.Range("K134").Value = 7874.823600000001
$ws.Range("L134").Value = 6214.2855
$ws.Range("M134").Value = -5339.823600000001
$ws.Range("N134").Value = -11284.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 8889128
$ws.Range("I70").Value = 8889128
$ws.Range("K70").Value = 26667384
$ws.Range("M70").Value = -26667069
$ws.Range("H73").Value = 8889128
$ws.Range("I73").Value = 8889128
$ws.Range("K73").Value = 26667384
$ws.Range("M73").Value = -26666292
$ws.Range("H75").Value = 200004350
$ws.Range("I75").Value = 800
$ws.Range("J75").Value = 250005250
$ws.Range("K75").Value = 2400
$ws.Range("L75").Value = 750015750
$ws.Range("M75").Value = -1402
$ws.Range("N75").Value = -750017746
$ws.Range("H78").Value = 200004350
$ws.Range("I78").Value = 800
$ws.Range("J78").Value = 250005250
$ws.Range("K78").Value = 7200
$ws.Range("L78").Value = 2250047250
$ws.Range("M78").Value = -2208
$ws.Range("N78").Value = -2250057234
$ws.Range("H131").Value = 2687.182
$ws.Range("I131").Value = 7400
$ws.Range("J131").Value = 1999.8959
$ws.Range("K131").Value = 22200
$ws.Range("L131").Value = 5999.6877
$ws.Range("M131").Value = -17160
$ws.Range("N131").Value = -16079.6877

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8483
$ws.Range("J92").Value = 8483
$ws.Range("L92").Value = 8483
$ws.Range("N92").Value = -12227

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1946.125
$ws.Range("I7").Value = 1674.8
$ws.Range("J7").Value = 2398.3333
$ws.Range("K7").Value = 1674.8
$ws.Range("L7").Value = 2398.3333
$ws.Range("M7").Value = -1562.8
$ws.Range("N7").Value = -2622.3333
$ws.Range("H40").Value = 2130.75
$ws.Range("I40").Value = 1790.2858
$ws.Range("K40").Value = 1790.2858
$ws.Range("M40").Value = -1654.2858
$ws.Range("H100").Value = 2412.7856
$ws.Range("I100").Value = 1734.875
$ws.Range("J100").Value = 3316.6667
$ws.Range("K100").Value = 1734.875
$ws.Range("L100").Value = 3316.6667
$ws.Range("M100").Value = -1193.875
$ws.Range("N100").Value = -4398.6667
$ws.Range("H126").Value = 1946.125
$ws.Range("I126").Value = 1674.8
$ws.Range("J126").Value = 2398.3333
$ws.Range("K126").Value = 5024.4
$ws.Range("L126").Value = 7194.999899999999
$ws.Range("M126").Value = -2554.4
$ws.Range("N126").Value = -12134.9999
